$d = $word.ActiveDocument
$paras = $d.Paragraphs

# The "KEY ACHIEVEMENTS AND IMPACT" section contains an "Impact" sub-heading
# followed by six bullet paragraphs, then the "TECHNICAL SKILLS" heading.
# Some of this bullet text also appears verbatim elsewhere in the resume
# (e.g. under "Partner - Siege Analytics"), so all Find/Replace and delete
# operations below are scoped to a Range between those two headings.

$headingStart = $null
$headingEnd = $null
$i = 0
foreach ($p in $paras) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t.Trim() -eq "Impact" -and $headingStart -eq $null) {
        $headingStart = $i
    }
    if ($t -like "*TECHNICAL SKILLS*" -and $headingStart -ne $null -and $headingEnd -eq $null) {
        $headingEnd = $i
    }
}

$scopeStartPos = $paras.Item($headingStart).Range.End

# Bookmark the start of the "TECHNICAL SKILLS" paragraph so the end of our
# working range auto-adjusts as the bullets above it are edited/removed.
$d.Bookmarks.Add("zzScopeEnd", $paras.Item($headingEnd).Range) | Out-Null

function Get-ScopedRange {
    $endPos = $d.Bookmarks.Item("zzScopeEnd").Range.Start
    return $d.Range($scopeStartPos, $endPos)
}

# Rewrite bullets 1-3 as impact-focused accomplishment statements.
$s = Get-ScopedRange
$s.Find.Execute(
    "Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%",
    2) | Out-Null

$s = Get-ScopedRange
$s.Find.Execute(
    "Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "`$4.7M savings enabled nonprofit access",
    2) | Out-Null

$s = Get-ScopedRange
$s.Find.Execute(
    "Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions",
    2) | Out-Null

# Rewrite bullet 6 (the last bullet).
$s = Get-ScopedRange
$s.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "178% accuracy improvement in racial classification algorithms",
    2) | Out-Null

# Remove bullets 4 and 5 entirely (their content was dropped, not rewritten).
# Use Find (no replacement) to locate each bullet, then Expand(4) =
# wdParagraph so the deleted range includes the full paragraph (and its
# paragraph mark), not just the matched text.
$s = Get-ScopedRange
$s.Find.Execute(
    "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis") | Out-Null
$s.Expand(4) | Out-Null
$bullet4Start = $s.Start

$s2 = Get-ScopedRange
$s2.Find.Execute(
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations") | Out-Null
$s2.Expand(4) | Out-Null
$bullet5End = $s2.End

$d.Range($bullet4Start, $bullet5End).Delete() | Out-Null

$d.Bookmarks.Item("zzScopeEnd").Delete() | Out-Null
